$wb = $excel.ActiveWorkbook

# 1. Rename the "CM" sheet to "AM" (position / sheetId / rId unchanged)
$wsAM = $wb.Worksheets.Item("CM")
$wsAM.Name = "AM"

# 2. Rename the table column headers (row 6) on each position sheet.
#    Renaming the header cell text updates both the shared-string table and
#    the underlying ListObject/table column name automatically.

# CB sheet
$wsCB = $wb.Worksheets.Item("CB")
$wsCB.Range("I6").Value = "Skilled Center Back"
$wsCB.Range("J6").Value = "Sweeper Center Back"
$wsCB.Range("K6").Value = "Imperative Center Back"

# FB sheet
$wsFB = $wb.Worksheets.Item("FB")
$wsFB.Range("I6").Value = "Defensive Full Back"
$wsFB.Range("J6").Value = "Offensive Full Back"
$wsFB.Range("K6").Value = "Imperative Full Back"

# DM sheet
$wsDM = $wb.Worksheets.Item("DM")
$wsDM.Range("I6").Value = "Creative Central Midfielder"
$wsDM.Range("J6").Value = "Defensive Central Midfielder"

# AM sheet (formerly CM)
$wsAM.Range("I6").Value = "Box to Box Attacking Midfielder"
$wsAM.Range("J6").Value = "Builder Attacking Midfielder"
$wsAM.Range("K6").Value = "Space Invader Attacking Midfielder"

# 3. Update the remembered selections on each sheet to match where the
#    author last left the cursor before saving.
[void]$wsCB.Range("K7").Select()
[void]$wsDM.Range("I6").Select()
[void]$wsAM.Range("K6").Select()

$wsW = $wb.Worksheets.Item("W")
[void]$wsW.Range("J6").Select()

# 4. Finally, the workbook was saved with the CF tab active/selected.
$wsCF = $wb.Worksheets.Item("CF")
[void]$wsCF.Activate()
